$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 into a new row 3, preserving its formatting exactly.
$ws.Rows("2:2").Copy() | Out-Null
$ws.Rows("3:3").Insert(-4121) | Out-Null

# This new row represents another (non-HMIS) client with no assessments,
# so only the Homebase ID differs from the copied template row.
$ws.Range("D3").Value = 24212

$ws.Range("D4").Select() | Out-Null
